# Regenerate column G ("K") values for the saved data sheet.
# The underlying model data (std/mean, s_vals) that produces these K values
# was recalculated upstream; here we write the resulting values directly
# into the worksheet, replacing the previous "Strike#" derived numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New K values for rows 2..39 (row 1 is the header row)
$kValues = @(5,5,6,9,7,8,7,6,9,5,5,5,6,10,13,6,8,3,3,8,9,7,8,5,7,7,5,5,5,3,8,3,3,5,2,5,0,2)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
